$wb = $excel.ActiveWorkbook

# --- Overview sheet: update Status cells for the ede2707b row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: Status for row 3, Error Detail (column P) for row 3 & widen column P ---
# (the saved col width ends up 5/6 wider than the ColumnWidth we assign, so
# back that off here to land on an on-disk width of exactly 40)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"
$wsZh.Range("P3").Value = "Handback file name: ffbk1gnc.ad1 is different with handoff file name: ede2707b-4988-4908-9089-7114a0198601.f240af75bafe8c0024344a2dca10141d5e9c7038.zh-cn."
$wsZh.Range("P1:P3").ColumnWidth = 39.16666666666667

# --- de-de sheet: Status for row 3, Error Detail (column P) for row 3 & widen column P ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"
$wsDe.Range("P3").Value = "Handback file name: ffbk1gnc.ad1 is different with handoff file name: ede2707b-4988-4908-9089-7114a0198601.f240af75bafe8c0024344a2dca10141d5e9c7038.de-de."
$wsDe.Range("P1:P3").ColumnWidth = 39.16666666666667
